$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new question row (row 23) that the diff introduces.
$ws.Range("A23").Value = "What is the primary use of OpenCV?"
$ws.Range("B23").Value = "A library for numerical computations in Python.
A library for image processing and computer vision.
A framework for developing web applications.
A tool for machine learning and deep learning."
# The "Answer Index" for this row is stored as text ("1"), not a number,
# so prefix with an apostrophe to force text entry like Excel would.
$ws.Range("C23").Value = "'1"
$ws.Range("D23").Value = "easy"

# Writing the multi-line text into B23 makes Excel auto-expand the row
# height; re-run AutoFit so the row settles back to the default height
# (matching the target workbook, which has no custom row height here).
$ws.Rows(23).EntireRow.AutoFit()
